$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 updates
$ws.Range("G6").Value = 2.8
$ws.Range("I6").Value = 2.7
$ws.Range("J6").Value = 3.6
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 11
$ws.Range("AA6").Value = 26
$ws.Range("AI6").Value = 12
$ws.Range("AO6").Value = 17
$ws.Range("AQ6").Value = 51
$ws.Range("AS6").Value = 301
$ws.Range("AV6").Value = 67
$ws.Range("AY6").Value = 29

# Row 7 updates
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.67

# Row 11 updates
$ws.Range("H11").Value = 3.45
$ws.Range("J11").Value = 2.07
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 6.3
$ws.Range("O11").Value = 1.34
$ws.Range("P11").Value = 2.75
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.65
$ws.Range("T11").Value = 2.52
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.65
$ws.Range("W11").Value = 5.6
$ws.Range("X11").Value = 6.5
$ws.Range("Y11").Value = 8
$ws.Range("AA11").Value = 13.5
$ws.Range("AC11").Value = 8
$ws.Range("AE11").Value = 19
$ws.Range("AH11").Value = 14.5
$ws.Range("AJ11").Value = 21
$ws.Range("AL11").Value = 90
$ws.Range("AM11").Value = 80
$ws.Range("AO11").Value = 7.3
$ws.Range("AP11").Value = 17.5
$ws.Range("AQ11").Value = 23
$ws.Range("AR11").Value = 55
$ws.Range("AU11").Value = 7.8
$ws.Range("AV11").Value = 80
$ws.Range("AW11").Value = 7.7
$ws.Range("AX11").Value = 40
$ws.Range("AZ11").Value = 300
$ws.Range("BA11").Value = 300
